# Re-order the fantasy roster table (rows 2-19) on the active sheet so each
# player/position/team trio ends up on the row dictated by the new layout
# (e.g. "Jalen Brunson" moves up to row 2, "Coby White" down to row 3, ...).
# The player names, their positions and teams are written out explicitly
# for every row so the whole table ends up in the exact target order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jalen Brunson",     "PG",       "New York Knicks"),
    @("Coby White",        "PG,SG",    "Chicago Bulls"),
    @("Trae Young",        "PG",       "Atlanta Hawks"),
    @("Jamal Murray",      "PG,SG",    "Denver Nuggets"),
    @("LeBron James",      "SF,PF",    "Los Angeles Lakers"),
    @("P.J. Washington",   "PF",       "Dallas Mavericks"),
    @("Kawhi Leonard",     "SG,SF,PF", "LA Clippers"),
    @("Devin Vassell",     "SG,SF",    "San Antonio Spurs"),
    @("Clint Capela",      "C",        "Atlanta Hawks"),
    @("Walker Kessler",    "C",        "Utah Jazz"),
    @("Desmond Bane",      "SG,SF",    "Memphis Grizzlies"),
    @("Norman Powell",     "SG,SF",    "LA Clippers"),
    @("Devin Booker",      "PG,SG",    "Phoenix Suns"),
    @("Myles Turner",      "C",        "Indiana Pacers"),
    @("Luguentz Dort",     "SG,SF",    "Oklahoma City Thunder"),
    @("D'Angelo Russell",  "PG",       "Brooklyn Nets"),
    @("Brandon Ingram",    "SG,SF,PF", "New Orleans Pelicans"),
    @("Immanuel Quickley", "PG,SG",    "Toronto Raptors")
)

$row = 2
foreach ($entry in $data) {
    $name = $entry[0]
    $pos = $entry[1]
    $team = $entry[2]

    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $pos
    $ws.Cells.Item($row, 3).Value = $team

    $row = $row + 1
}
